# Daily attendance processing - 2025-11-09 18:49:53
# Adds 24 new ANATOMY Session 2 check-ins to the Attendance sheet
# and refreshes the corresponding rows in the Summary sheet.

$wb = $excel.ActiveWorkbook
$wsAtt = $wb.Worksheets.Item("Attendance")
$wsSum = $wb.Worksheets.Item("Summary")

# ---- Attendance sheet: append new check-in rows (371-394) ----
$newRows = @(
    @{ Row=371; A='221966'; B='ساكى جوزيف اليساما زونقبيتى'; C='Year 2'; D='C1'; E='221966@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:26:29'; K='C1' }
    @{ Row=372; A='221599'; B='سلمى عبد الرحمن عبيد موسى'; C='Year 2'; D='C1'; E='221599@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:26:33'; K='C1' }
    @{ Row=373; A='221909'; B='ديكتور يمبيك بول نيان'; C='Year 2'; D='C1'; E='221909@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:26:39'; K='C1' }
    @{ Row=374; A='221774'; B='يدجوك جيمس كوانقو اكوك'; C='Year 2'; D='C1'; E='221774@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:26:41'; K='C1' }
    @{ Row=375; A='221938'; B='مانويلا ناكوتا مارينو لوكالى'; C='Year 2'; D='C1'; E='221938@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:26:41'; K='C1' }
    @{ Row=376; A='221996'; B='نياقوط فال توت دوير'; C='Year 2'; D='C1'; E='221996@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:26:43'; K='C1' }
    @{ Row=377; A='211620'; B='محمودول اسلام'; C='Year 2'; D='C1'; E='211620@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:26:48'; K='C1' }
    @{ Row=378; A='222053'; B='صباح سيف الدين عثمان اسحق'; C='Year 2'; D='C1'; E='222053@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:26:56'; K='C1' }
    @{ Row=379; A='211776'; B='رقيه ادريسو'; C='Year 2'; D='C1'; E='211776@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:28:58'; K='C1' }
    @{ Row=380; A='221822'; B='سعادة يوسف عليو'; C='Year 2'; D='C1'; E='221822@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:44:55'; K='C1' }
    @{ Row=381; A='221866'; B='امينة موسى رمبو'; C='Year 2'; D='C1'; E='221866@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:45:01'; K='C1' }
    @{ Row=382; A='221914'; B='معز اشتياق'; C='Year 2'; D='C1'; E='221914@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:45:06'; K='C1' }
    @{ Row=383; A='221755'; B='سعدية عاشق'; C='Year 2'; D='C1'; E='221755@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:45:12'; K='C1' }
    @{ Row=384; A='221539'; B='تحريم شوكات مالك'; C='Year 2'; D='C1'; E='221539@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:45:16'; K='C1' }
    @{ Row=385; A='222032'; B='احمد شوقى عبد الرحيم طه'; C='Year 2'; D='C1'; E='222032@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:45:22'; K='C1' }
    @{ Row=386; A='222002'; B='يوسف عثمان باباغانا'; C='Year 2'; D='C1'; E='222002@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:45:30'; K='C1' }
    @{ Row=387; A='222027'; B='خيرية عبد الرازق'; C='Year 2'; D='C1'; E='222027@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:45:36'; K='C1' }
    @{ Row=388; A='212386'; B='لويد اليكس موجا'; C='Year 2'; D='C1'; E='212386@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:45:49'; K='C1' }
    @{ Row=389; A='221695'; B='اروب ميجوك دودى باقوير'; C='Year 2'; D='C1'; E='221695@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:46:10'; K='C1' }
    @{ Row=390; A='221329'; B='خديجة اولو اتوين ادونبكو'; C='Year 2'; D='C1'; E='221329@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:46:18'; K='C1' }
    @{ Row=391; A='221826'; B='بخيتة اوت قور كول'; C='Year 2'; D='C1'; E='221826@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:46:25'; K='C1' }
    @{ Row=392; A='221982'; B='امنويل اكوى اقوتو كوت'; C='Year 2'; D='C1'; E='221982@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:46:34'; K='C1' }
    @{ Row=393; A='221923'; B='ليفو سوزى وليام جوزيف'; C='Year 2'; D='C1'; E='221923@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:46:59'; K='C1' }
    @{ Row=394; A='212145'; B='سميه لاوان شايبو'; C='Year 2'; D='C1'; E='212145@med.asu.edu.eg'; F='ANATOMY'; G='2'; H='ANATOMY'; I='09/11/2025'; J='14:47:35'; K='C1' }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K")) {
        $cell = $wsAtt.Range("$col$r")
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$col]
        $cell.ClearFormats()
    }
}

# Refresh the autofilter range to cover the newly appended rows
$wsAtt.AutoFilterMode = $false
$wsAtt.Range("A1:K394").AutoFilter()

# Update the hidden _xlnm._FilterDatabase defined name for the Attendance sheet
for ($i = 1; $i -le $wb.Names.Count(); $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name() -eq "Attendance!_FilterDatabase") {
        $nm.RefersTo = "='Attendance'!`$A`$1:`$K`$394"
    }
}

# ---- Summary sheet: recompute risk/attendance stats for affected students ----
$riskColor = @{
    "High Risk"     = 8158463
    "Moderate Risk" = 8239615
    "Low Risk"      = 10940927
}

$summaryUpdates = @(
    @{ Row=37; F='Moderate Risk'; G='3.4%'; H=21; L=1; M=5; O=1; Q=1 }
    @{ Row=39; F=$null; G='6.9%'; H=20; L=2; M=4; O=1; Q=1 }
    @{ Row=52; F=$null; G='6.9%'; H=20; L=2; M=4; O=1; Q=1 }
    @{ Row=65; F=$null; G='13.8%'; H=18; L=4; M=2; O=2; Q=1 }
    @{ Row=112; F='Moderate Risk'; G='3.4%'; H=21; L=1; M=5; O=1; Q=1 }
    @{ Row=147; F='Low Risk'; G='10.3%'; H=19; L=3; M=3; O=1; Q=1 }
    @{ Row=160; F='Low Risk'; G='10.3%'; H=19; L=3; M=3; O=2; Q=1 }
    @{ Row=178; F='Moderate Risk'; G='3.4%'; H=21; L=1; M=5; O=1; Q=1 }
    @{ Row=184; F='Low Risk'; G='10.3%'; H=19; L=3; M=3; O=1; Q=1 }
    @{ Row=188; F=$null; G='6.9%'; H=20; L=2; M=4; O=1; Q=1 }
    @{ Row=195; F='Low Risk'; G='10.3%'; H=19; L=3; M=3; O=1; Q=1 }
    @{ Row=196; F='Low Risk'; G='10.3%'; H=19; L=3; M=3; O=1; Q=1 }
    @{ Row=203; F=$null; G='6.9%'; H=20; L=2; M=4; O=1; Q=1 }
    @{ Row=210; F='Low Risk'; G='10.3%'; H=19; L=3; M=3; O=1; Q=1 }
    @{ Row=212; F=$null; G='6.9%'; H=20; L=2; M=4; O=1; Q=1 }
    @{ Row=214; F=$null; G='6.9%'; H=20; L=2; M=4; O=1; Q=1 }
    @{ Row=217; F=$null; G='13.8%'; H=18; L=4; M=2; O=1; Q=1 }
    @{ Row=225; F='Low Risk'; G='10.3%'; H=19; L=3; M=3; O=1; Q=1 }
    @{ Row=227; F=$null; G='13.8%'; H=18; L=4; M=2; O=2; Q=1 }
    @{ Row=228; F=$null; G='6.9%'; H=20; L=2; M=4; O=1; Q=1 }
    @{ Row=231; F=$null; G='6.9%'; H=20; L=2; M=4; O=1; Q=1 }
    @{ Row=237; F=$null; G='6.9%'; H=20; L=2; M=4; O=1; Q=1 }
    @{ Row=240; F='Low Risk'; G='10.3%'; H=19; L=3; M=3; O=1; Q=1 }
    @{ Row=243; F=$null; G='13.8%'; H=18; L=4; M=2; O=2; Q=1 }
)

foreach ($u in $summaryUpdates) {
    $r = $u.Row
    if ($u.F) {
        $wsSum.Range("F$r").Value = $u.F
        $wsSum.Range("F$r").Interior.Color = $riskColor[$u.F]
        $wsSum.Range("F$r").Font.Bold = $true
    }
    $wsSum.Range("G$r").NumberFormat = "@"
    $wsSum.Range("G$r").Value = $u.G
    $wsSum.Range("G$r").NumberFormat = "0.0%"
    $wsSum.Range("H$r").Value = $u.H
    $wsSum.Range("L$r").Value = $u.L
    $wsSum.Range("M$r").Value = $u.M
    $wsSum.Range("O$r").Value = $u.O
    $wsSum.Range("Q$r").Value = $u.Q
}

Write-Host "Attendance rows added:" $newRows.Count
Write-Host "Summary rows updated:" $summaryUpdates.Count
